# Update the cryptocurrency price / 1h-volume table to the latest
# scraped snapshot (GitHub Actions refresh), and fix the row 44/45
# coin ordering (dogwifhat now ranks above USDe).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.897.46"
$ws.Range("E2").Value = "'  +2.75%  "
$ws.Range("D3").Value = "'2.529.48"
$ws.Range("E3").Value = "'  +1.36%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'592.80"
$ws.Range("E5").Value = "'  +2.32%  "
$ws.Range("D6").Value = "'177.06"
$ws.Range("E6").Value = "'  +4.29%  "
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("D8").Value = "'0.532"
$ws.Range("E8").Value = "'  +2.53%  "
$ws.Range("D9").Value = "'2.527.23"
$ws.Range("E9").Value = "'  +1.31%  "
$ws.Range("E10").Value = "'  +3.23%  "
$ws.Range("E11").Value = "'  +3.15%  "
$ws.Range("E12").Value = "'  +1.49%  "
$ws.Range("E13").Value = "'  -0.39%  "
$ws.Range("D14").Value = "'26.90"
$ws.Range("E14").Value = "'  +2.18%  "
$ws.Range("D15").Value = "'2.987.42"
$ws.Range("E15").Value = "'  +1.63%  "
$ws.Range("E16").Value = "'  +2.69%  "
$ws.Range("D17").Value = "'67.650.59"
$ws.Range("E17").Value = "'  +2.83%  "
$ws.Range("D18").Value = "'2.525.04"
$ws.Range("E18").Value = "'  +0.91%  "
$ws.Range("D19").Value = "'8.04"
$ws.Range("E19").Value = "'  +5.44%  "
$ws.Range("D20").Value = "'11.46"
$ws.Range("E20").Value = "'  +2.78%  "
$ws.Range("D21").Value = "'364.56"
$ws.Range("E21").Value = "'  +6.19%  "
$ws.Range("D22").Value = "'4.20"
$ws.Range("E22").Value = "'  +0.80%  "
$ws.Range("D23").Value = "'4.67"
$ws.Range("E23").Value = "'  +2.90%  "
$ws.Range("D24").Value = "'1.96"
$ws.Range("E24").Value = "'  +1.41%  "
$ws.Range("E25").Value = "'  -0.07%  "
$ws.Range("D26").Value = "'71.12"
$ws.Range("E26").Value = "'  +3.03%  "
$ws.Range("E27").Value = "'  +4.37%  "
$ws.Range("E28").Value = "'  -0.20%  "
$ws.Range("D29").Value = "'2.658.06"
$ws.Range("E29").Value = "'  +1.46%  "
$ws.Range("D30").Value = "'0.0₃0991"
$ws.Range("E30").Value = "'  +3.37%  "
$ws.Range("D31").Value = "'545.69"
$ws.Range("E31").Value = "'  +4.29%  "
$ws.Range("D32").Value = "'8.30"
$ws.Range("E32").Value = "'  +3.42%  "
$ws.Range("E33").Value = "'  +2.56%  "
$ws.Range("D34").Value = "'1.87"
$ws.Range("E34").Value = "'  +3.16%  "
$ws.Range("E35").Value = "'  -0.80%  "
$ws.Range("E36").Value = "'  +0.01%  "
$ws.Range("E37").Value = "'  +1.64%  "
$ws.Range("D38").Value = "'156.43"
$ws.Range("E38").Value = "'  +0.06%  "
$ws.Range("D39").Value = "'18.85"
$ws.Range("E39").Value = "'  +2.18%  "
$ws.Range("D40").Value = "'18.67"
$ws.Range("E41").Value = "'  +1.61%  "
$ws.Range("D42").Value = "'5.20"
$ws.Range("E42").Value = "'  +3.22%  "
$ws.Range("D43").Value = "'1.80"
$ws.Range("E43").Value = "'  +1.97%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.52"
$ws.Range("E44").Value = "'  +4.19%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "'  -0.05%  "
$ws.Range("D46").Value = "'0.563"
$ws.Range("E46").Value = "'  +1.85%  "
$ws.Range("D47").Value = "'146.56"
$ws.Range("E47").Value = "'  +0.32%  "
$ws.Range("D48").Value = "'3.73"
$ws.Range("E48").Value = "'  +2.18%  "
$ws.Range("E49").Value = "'  +2.65%  "
$ws.Range("D51").Value = "'0.0757"
$ws.Range("E51").Value = "'  +0.97%  "
